# Add files via upload
#
# Slide 6 ("STAMP - Session-Sender Control Code Field"): reposition/resize
# the title placeholder and append " - Usage" to its title text.
#
# Slide 8 ("STAMP - Return Path TLV"): append " - Usage" to its title text.

$p = $ppt.ActivePresentation

$slide6 = $p.Slides.Item(6)
$title6 = $slide6.Shapes.Item(1)

$title6.Left = 5254 / 914400 * 72
$title6.Top = 103031 / 914400 * 72
$title6.Width = 9062545 / 914400 * 72
$title6.Height = 857250 / 914400 * 72

$title6.TextFrame.TextRange.Text = "STAMP - Session-Sender Control Code Field - Usage"

$slide8 = $p.Slides.Item(8)
$title8 = $slide8.Shapes.Item(1)

$title8.TextFrame.TextRange.Text = "STAMP - Return Path TLV - Usage"
